$d = $word.ActiveDocument

$replacements = @(
    @("83×88=7304", "75×54=4050"),
    @("85×93=7905", "97×40=3880"),
    @("74×88=6512", "57×94=5358"),
    @("54×18=972",  "36×60=2160"),
    @("51×85=4335", "30×60=1800"),
    @("23×32=736",  "66×36=2376"),
    @("91×25=2275", "94×32=3008"),
    @("26×28=728",  "93×82=7626"),
    @("13×21=273",  "97×21=2037"),
    @("79×14=1106", "64×61=3904"),
    @("50×85=4250", "92×61=5612"),
    @("55×99=5445", "42×62=2604"),
    @("20×88=1760", "21×25=525"),
    @("25×60=1500", "51×43=2193"),
    @("34×55=1870", "79×57=4503"),
    @("94×33=3102", "47×29=1363"),
    @("48×60=2880", "36×28=1008"),
    @("15×74=1110", "22×77=1694"),
    @("80×89=7120", "61×38=2318"),
    @("20×16=320",  "62×16=992"),
    @("73×41=2993", "13×32=416"),
    @("46×82=3772", "67×44=2948"),
    @("93×64=5952", "11×57=627"),
    @("26×22=572",  "11×67=737"),
    @("29×95=2755", "27×79=2133")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
